$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2024-08-12 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-08-13 Tuesday", 2)

# Update the multiplication problems in the table, cell by cell to avoid
# ambiguity between old/new values that collide across cells.
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="82×99="; New="50×12="},
    @{Row=1;  Col=2; Old="96×31="; New="27×42="},
    @{Row=1;  Col=3; Old="27×73="; New="27×32="},
    @{Row=1;  Col=4; Old="99×94="; New="28×21="},
    @{Row=1;  Col=5; Old="43×20="; New="63×38="},

    @{Row=5;  Col=1; Old="14×35="; New="60×45="},
    @{Row=5;  Col=2; Old="60×53="; New="49×79="},
    @{Row=5;  Col=3; Old="50×58="; New="34×40="},
    @{Row=5;  Col=4; Old="20×78="; New="30×49="},
    @{Row=5;  Col=5; Old="45×78="; New="74×47="},

    @{Row=10; Col=1; Old="53×20="; New="63×48="},
    @{Row=10; Col=2; Old="55×74="; New="14×35="},
    @{Row=10; Col=3; Old="28×56="; New="87×31="},
    @{Row=10; Col=4; Old="32×49="; New="29×41="},
    @{Row=10; Col=5; Old="50×60="; New="25×54="},

    @{Row=15; Col=1; Old="46×65="; New="67×32="},
    @{Row=15; Col=2; Old="53×48="; New="74×63="},
    @{Row=15; Col=3; Old="56×57="; New="71×38="},
    @{Row=15; Col=4; Old="23×33="; New="95×46="},
    @{Row=15; Col=5; Old="25×90="; New="17×97="},

    @{Row=20; Col=1; Old="37×83="; New="22×76="},
    @{Row=20; Col=2; Old="74×43="; New="35×50="},
    @{Row=20; Col=3; Old="32×64="; New="48×78="},
    @{Row=20; Col=4; Old="50×40="; New="40×41="},
    @{Row=20; Col=5; Old="27×68="; New="91×80="}
)

foreach ($rep in $replacements) {
    $cell = $t.Cell($rep.Row, $rep.Col)
    $cell.Range.Find.Execute($rep.Old, $true, $false, $false, $false, $false,
                              $true, 1, $false, $rep.New, 2)
}
